$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the per-row "=Fn*En" formulas in G3:G25 into one shared-formula
# group (mirrors what Excel does when the same formula is (re)entered across
# a contiguous range), matching the shared formula that already exists for
# G26.
$ws.Range("G3:G25").Formula = "=F3*E3"

# Home Depot wood price correction (row 26): 20.00 -> 19.63. This ripples
# into G26, the subtotal (G28) and the grand total (G30) automatically.
$ws.Range("E26").Value = 19.63

# Restore the view: no saved scroll anchor, selection resting on the
# subtotal cell (merged G28:H28).
$ws.Range("G28:H28").Select() | Out-Null
